$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7070478796958923
$ws.Range("B1").Value = 1.185357928276062
$ws.Range("C1").Value = 1.645100951194763
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 15
